# Update the Malta MSME country-indicator figures in the "Summary" sheet.
# These figures are stored as text (numbers-as-strings) in the workbook, so
# force the target cells to a text number format before writing the new
# values — this keeps them stored as text (matching the original data)
# instead of being auto-converted to floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$targets = @("B13", "C13", "D13", "B14", "C14", "D14", "B16", "C16", "D16")
foreach ($addr in $targets) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 13: Enterprises density (per 1000 people) -- Micro / SMEs / MSMEs
$ws.Range("B13").Value = "59.37"
$ws.Range("C13").Value = "3.61"
$ws.Range("D13").Value = "62.97"

# Row 14: Employment (% of total) -- Micro / SMEs / MSMEs
$ws.Range("B14").Value = "35.97"
$ws.Range("C14").Value = "44.07"
$ws.Range("D14").Value = "80.03"

# Row 16: Enterprises (% of total) -- Micro / SMEs / MSMEs
$ws.Range("B16").Value = "94.11"
$ws.Range("C16").Value = "5.72"
$ws.Range("D16").Value = "99.83"
